# Added CO2 output column and fixed the max/ordering of the tkm-*Usage
# columns at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("crudeoil") to make room for the
# new "CO2" output. Everything from C onward (through AI) shifts one
# column to the right (D..AJ).
$ws.Columns("C:C").Insert()

# Populate the newly inserted CO2 column.
$ws.Range("C1").Value = "CO2"
$ws.Range("C2").Value = 87.14279999999999

# Fix the ordering of the last five "Usage" columns (previously
# tkm-N2Usage, tkm-N3Usage, tkm-N1Usage, pkmUsage, tkm-SZMUsage - now
# shifted to AF:AJ by the insert above). The corrected order moves
# tkm-SZMUsage to the front of this group and tkm-N3Usage to the back,
# while tkm-N1Usage / pkmUsage keep their relative order; the
# underlying values for each named column are unchanged.
$ws.Range("AF1").Value = "tkm-SZMUsage"
$ws.Range("AG1").Value = "tkm-N2Usage"
$ws.Range("AH1").Value = "tkm-N1Usage"
$ws.Range("AI1").Value = "pkmUsage"
$ws.Range("AJ1").Value = "tkm-N3Usage"

$ws.Range("AF2").Value = 414.5
$ws.Range("AG2").Value = 24.2
$ws.Range("AH2").Value = 7.5
$ws.Range("AI2").Value = 858
$ws.Range("AJ2").Value = 130.3
